# INS03 Segunda Aba P Média
# - Renames the single sheet "INS03 VIP 22,5M" -> "INS03 Prioridade Baixa"
# - Reworks the hidden helper column Y into a visible "Sla Consumido" column
#   and adds a new "Sim/Não" column Z, with the B4 KPI formula now reading
#   column Z with a "=0" test instead of column Y with a "<=1350" test.
# - Duplicates the (now updated) sheet into a second tab
#   "INS03 Prioridade Média" whose only content difference is the A4 label.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# ---- Sheet 1: "INS03 VIP 22,5M" -> "INS03 Prioridade Baixa" -------------
$ws1.Name = "INS03 Prioridade Baixa"

# Column Y: was hidden (helper column) -> now shown, same width/style as X
$ws1.Columns(25).Hidden = $false
$ws1.Columns(25).ColumnWidth = $ws1.Columns(24).ColumnWidth()

# Column Z: brand new helper column
$ws1.Columns(26).ColumnWidth = $ws1.Columns(6).ColumnWidth()

# KPI label + formula
$ws1.Range("A4").Value = "Total de Tickets Solucionados Prioridade Baixa"
$ws1.Range("B4").Formula = '=COUNTIF(Z9:Z9,"=0")'

# New header cells on row 8: Y8 "Sla Consumido", Z8 "Sim/Não"
# (set the text first, then pull the Q8:X8-style formatting via PasteSpecial
# so the new cells land on the same style as the rest of the header row)
$ws1.Range("Y8").Value = "Sla Consumido"
$ws1.Range("X8").Copy()
$ws1.Range("Y8").PasteSpecial(-4122)

$ws1.Range("Z8").Value = "Sim/Não"
$ws1.Range("X8").Copy()
$ws1.Range("Z8").PasteSpecial(-4122)

# ---- Sheet 2: duplicate of sheet 1, "INS03 Prioridade Média" ------------
$ws1.Copy($null, $ws1)
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "INS03 Prioridade Média"
$ws2.Range("A4").Value = "Total de Tickets Solucionados Prioridade Média"

# Keep "INS03 Prioridade Baixa" (sheet 1) as the selected / active tab,
# same as the original workbook.
$ws1.Select()
